# Matrix.xlsx - "Added file read feature"
# Update the Sheet1 matrix values and move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = 10
$ws.Range("C1").Value = 13
$ws.Range("A2").Value = 14

$ws.Range("A2").Select()
